$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.416.39"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "3.050.47"
$ws.Range("E3").Value = "  +4.47%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.209"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.84%  "
$ws.Range("D10").Value = "3.046.77"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  +5.32%  "
$ws.Range("D14").Value = "3.611.36"
$ws.Range("E14").Value = "  +4.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").Value = "76.353.95"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000194"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").Value = "3.065.64"
$ws.Range("E18").Value = "  +5.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D25").Value = "3.206.66"
$ws.Range("E25").Value = "  +4.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +7.85%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "507.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +7.01%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.388"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.20%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "191.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.13%  "
$ws.Range("E41").Value = "  -6.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.785"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.612"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.10%  "
